$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = "line1`nline2`n"
$v = $ws.Range("A1").Value()
$v
